# Update TPM-derived LR-pair metrics for Ucn2-Crhr2 sheet (rows 2-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05692066666666667
$ws.Range("H2").Value = 0.170762
$ws.Range("I2").Value = 0.07235837399807114
$ws.Range("J2").Value = 0.07235837399807114
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2792176666666666
$ws.Range("N2").Value = 0.837653
$ws.Range("O2").Value = 0.06481516157886429
$ws.Range("P2").Value = 0.06481516157886429
$ws.Range("Q2").Value = 0.01589325573177778
$ws.Range("R2").Value = 0.143039301586
$ws.Range("S2").Value = 0.004689919702268873
$ws.Range("T2").Value = 0.004689919702268873
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.05692066666666667
$ws.Range("H3").Value = 0.170762
$ws.Range("I3").Value = 0.07235837399807114
$ws.Range("J3").Value = 0.07235837399807114
$ws.Range("O3").Value = 0.0127693110033334
$ws.Range("P3").Value = 0.0127693110033334
$ws.Range("Q3").Value = 0.003131148952666667
$ws.Range("R3").Value = 0.028180340574
$ws.Range("S3").Value = 0.0009239665812768836
$ws.Range("T3").Value = 0.0009239665812768836
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.05692066666666667
$ws.Range("H4").Value = 0.170762
$ws.Range("I4").Value = 0.07235837399807114
$ws.Range("J4").Value = 0.07235837399807114
$ws.Range("M4").Value = 3.97368
$ws.Range("N4").Value = 11.92104
$ws.Range("O4").Value = 0.9224155274178023
$ws.Range("P4").Value = 0.9224155274178023
$ws.Range("Q4").Value = 0.22618451472
$ws.Range("R4").Value = 2.03566063248
$ws.Range("S4").Value = 0.06674448771452539
$ws.Range("T4").Value = 0.06674448771452539
$ws.Range("I5").Value = 0.6056548703615503
$ws.Range("J5").Value = 0.6056548703615503
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2792176666666666
$ws.Range("N5").Value = 0.837653
$ws.Range("O5").Value = 0.06481516157886429
$ws.Range("P5").Value = 0.06481516157886429
$ws.Range("Q5").Value = 0.1330299066713333
$ws.Range("R5").Value = 1.197269160042
$ws.Range("S5").Value = 0.03925561828350999
$ws.Range("T5").Value = 0.03925561828350999
$ws.Range("I6").Value = 0.6056548703615503
$ws.Range("J6").Value = 0.6056548703615503
$ws.Range("O6").Value = 0.0127693110033334
$ws.Range("P6").Value = 0.0127693110033334
$ws.Range("S6").Value = 0.007733795400330212
$ws.Range("T6").Value = 0.007733795400330212
$ws.Range("I7").Value = 0.6056548703615503
$ws.Range("J7").Value = 0.6056548703615503
$ws.Range("M7").Value = 3.97368
$ws.Range("N7").Value = 11.92104
$ws.Range("O7").Value = 0.9224155274178023
$ws.Range("P7").Value = 0.9224155274178023
$ws.Range("Q7").Value = 1.89321215184
$ws.Range("R7").Value = 17.03890936656
$ws.Range("S7").Value = 0.5586654566777101
$ws.Range("T7").Value = 0.5586654566777101
$ws.Range("G8").Value = 0.2532906666666667
$ws.Range("H8").Value = 0.7598720000000001
$ws.Range("I8").Value = 0.3219867556403786
$ws.Range("J8").Value = 0.3219867556403785
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2792176666666666
$ws.Range("N8").Value = 0.837653
$ws.Range("O8").Value = 0.06481516157886429
$ws.Range("P8").Value = 0.06481516157886429
$ws.Range("Q8").Value = 0.07072322893511111
$ws.Range("R8").Value = 0.6365090604160001
$ws.Range("S8").Value = 0.02086962359308543
$ws.Range("T8").Value = 0.02086962359308543
$ws.Range("G9").Value = 0.2532906666666667
$ws.Range("H9").Value = 0.7598720000000001
$ws.Range("I9").Value = 0.3219867556403786
$ws.Range("J9").Value = 0.3219867556403785
$ws.Range("O9").Value = 0.0127693110033334
$ws.Range("P9").Value = 0.0127693110033334
$ws.Range("Q9").Value = 0.01393326628266667
$ws.Range("R9").Value = 0.125399396544
$ws.Range("S9").Value = 0.004111549021726311
$ws.Range("T9").Value = 0.00411154902172631
$ws.Range("G10").Value = 0.2532906666666667
$ws.Range("H10").Value = 0.7598720000000001
$ws.Range("I10").Value = 0.3219867556403786
$ws.Range("J10").Value = 0.3219867556403785
$ws.Range("M10").Value = 3.97368
$ws.Range("N10").Value = 11.92104
$ws.Range("O10").Value = 0.9224155274178023
$ws.Range("P10").Value = 0.9224155274178023
$ws.Range("Q10").Value = 1.00649605632
$ws.Range("R10").Value = 9.058464506880002
$ws.Range("S10").Value = 0.2970055830255668
$ws.Range("T10").Value = 0.2970055830255668
